$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.106.02"
$ws.Range("E2").Value = "  +0.18%  "
$ws.Range("D3").Value = "2.395.63"
$ws.Range("E3").Value = "  +4.70%  "
$ws.Range("E4").Value = "  -0.16%  "
$ws.Range("D5").Value = "'333.59"
$ws.Range("E5").Value = "  +7.51%  "
$ws.Range("D6").Value = "'104.35"
$ws.Range("E6").Value = "  -7.84%  "
$ws.Range("D7").Value = "'0.646"
$ws.Range("E7").Value = "  +2.07%  "
$ws.Range("E8").Value = "  +0.08%  "
$ws.Range("D9").Value = "'0.643"
$ws.Range("E9").Value = "  +4.68%  "
$ws.Range("D10").Value = "'41.81"
$ws.Range("E10").Value = "  -6.17%  "
$ws.Range("D11").Value = "'0.0936"
$ws.Range("E11").Value = "  +0.85%  "
$ws.Range("D12").Value = "'8.67"
$ws.Range("E12").Value = "  -1.94%  "
$ws.Range("D13").Value = "'1.05"
$ws.Range("E13").Value = "  -1.63%  "
$ws.Range("D14").Value = "'17.05"
$ws.Range("E14").Value = "  +9.96%  "
$ws.Range("E15").Value = "  +1.75%  "
$ws.Range("D16").Value = "2.759.73"
$ws.Range("E16").Value = "  +4.93%  "
$ws.Range("D17").Value = "2.403.20"
$ws.Range("E17").Value = "  +5.20%  "
$ws.Range("D18").Value = "43.136.59"
$ws.Range("E18").Value = "  +0.12%  "
$ws.Range("D19").Value = "'7.61"
$ws.Range("E19").Value = "  +5.52%  "
$ws.Range("D20").Value = "'0.0000108"
$ws.Range("E20").Value = "  +0.59%  "
$ws.Range("D21").Value = "'3.91"
$ws.Range("E21").Value = "  +7.35%  "
$ws.Range("D22").Value = "'77.10"
$ws.Range("E22").Value = "  +2.08%  "
$ws.Range("D23").Value = "'275.76"
$ws.Range("E23").Value = "  +7.33%  "
$ws.Range("D24").Value = "'2.39"
$ws.Range("E24").Value = "  -3.03%  "
$ws.Range("D25").Value = "'9.90"
$ws.Range("E25").Value = "  +10.28%  "
$ws.Range("D26").Value = "'11.89"
$ws.Range("E26").Value = "  +0.84%  "
$ws.Range("D27").Value = "'0.999"
$ws.Range("E27").Value = "  -0.06%  "
$ws.Range("D28").Value = "'24.24"
$ws.Range("E28").Value = "  +9.00%  "
$ws.Range("E29").Value = "  -1.87%  "
$ws.Range("D30").Value = "'175.16"
$ws.Range("E30").Value = "  -0.07%  "
$ws.Range("B31").Value = "WEMIXToken"
$ws.Range("C31").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D31").Value = "'3.16"
$ws.Range("E31").Value = "  -0.86%  "
$ws.Range("B32").Value = "InjectiveProtocol"
$ws.Range("C32").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D32").Value = "'36.77"
$ws.Range("E32").Value = "  -4.10%  "
$ws.Range("B33").Value = "Hedera"
$ws.Range("C33").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D33").Value = "'0.0931"
$ws.Range("E33").Value = "  +3.19%  "
$ws.Range("E34").Value = "  +6.26%  "
$ws.Range("E35").Value = "  +4.84%  "
$ws.Range("D36").Value = "'4.82"
$ws.Range("E36").Value = "  -4.10%  "
$ws.Range("D37").Value = "'4.06"
$ws.Range("E37").Value = "  -3.66%  "
$ws.Range("D38").Value = "'0.0366"
$ws.Range("E38").Value = "  -3.15%  "
$ws.Range("E39").Value = "  +3.53%  "
$ws.Range("E40").Value = "  +11.20%  "
$ws.Range("D41").Value = "'1.55"
$ws.Range("E41").Value = "  +12.25%  "
$ws.Range("D42").Value = "'0.235"
$ws.Range("E42").Value = "  +1.82%  "
$ws.Range("D43").Value = "'70.50"
$ws.Range("E43").Value = "  -3.49%  "
$ws.Range("B44").Value = "FirstDigitalUSD"
$ws.Range("C44").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D44").Value = "'1.00"
$ws.Range("E44").Value = "  +0.16%  "
$ws.Range("B45").Value = "Aave"
$ws.Range("C45").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D45").Value = "'119.65"
$ws.Range("E45").Value = "  +10.74%  "
$ws.Range("D46").Value = "'91.54"
$ws.Range("E46").Value = "  +43.08%  "
$ws.Range("D47").Value = "'12.24"
$ws.Range("E47").Value = "  -3.04%  "
$ws.Range("D48").Value = "'5.57"
$ws.Range("E48").Value = "  -2.72%  "
$ws.Range("B49").Value = "WOONetwork"
$ws.Range("C49").Value = "https://coinranking.com/coin/k-J3YwacF+woonetwork-woo"
$ws.Range("D49").Value = "'0.523"
$ws.Range("E49").Value = "  +18.98%  "
$ws.Range("B50").Value = "FraxShare"
$ws.Range("C50").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D50").Value = "'9.15"
$ws.Range("E50").Value = "  +3.48%  "
$ws.Range("D51").Value = "'1.31"
$ws.Range("E51").Value = "  +0.60%  "
